$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove last row (row 18, APOLLOHOSP) which no longer exists in the updated data
$ws.Rows.Item(18).Delete()

# Insert a new column H ("To Year High Profit") shifting old H (Total Year Up) to I
$ws.Columns.Item(8).Insert()

$ws.Cells.Item(1, 8).Value = 'To Year High Profit'

# Row 2: TRENT
$ws.Cells.Item(2, 1).Value = 'INE849A01020'
$ws.Cells.Item(2, 2).Value = 'TRENT'
$ws.Cells.Item(2, 3).Value = 'Speciality Retail'
$ws.Cells.Item(2, 4).Value = 7284.45
$ws.Cells.Item(2, 5).Value = 8345
$ws.Cells.Item(2, 6).Value = 2955
$ws.Cells.Item(2, 7).Value = 12.70880766926303
$ws.Cells.Item(2, 8).Value = 14.55909505865234
$ws.Cells.Item(2, 9).Value = 146.5126903553299

# Row 3: M&M
$ws.Cells.Item(3, 1).Value = 'INE101A01026'
$ws.Cells.Item(3, 2).Value = 'M&M'
$ws.Cells.Item(3, 3).Value = 'Passenger Cars & Utility Vehicles'
$ws.Cells.Item(3, 4).Value = 3180
$ws.Cells.Item(3, 5).Value = 3237.05
$ws.Cells.Item(3, 6).Value = 1575
$ws.Cells.Item(3, 7).Value = 1.762407129948573
$ws.Cells.Item(3, 8).Value = 1.794025157232704
$ws.Cells.Item(3, 9).Value = 101.9047619047619

# Row 4: BEL
$ws.Cells.Item(4, 1).Value = 'INE263A01024'
$ws.Cells.Item(4, 2).Value = 'BEL'
$ws.Cells.Item(4, 3).Value = 'Aerospace & Defense'
$ws.Cells.Item(4, 4).Value = 291.95
$ws.Cells.Item(4, 5).Value = 340.5
$ws.Cells.Item(4, 6).Value = 171.75
$ws.Cells.Item(4, 7).Value = 14.25844346549193
$ws.Cells.Item(4, 8).Value = 16.62955985613976
$ws.Cells.Item(4, 9).Value = 69.98544395924309

# Row 5: HCLTECH
$ws.Cells.Item(5, 1).Value = 'INE860A01027'
$ws.Cells.Item(5, 2).Value = 'HCLTECH'
$ws.Cells.Item(5, 3).Value = 'Computers - Software & Consulting'
$ws.Cells.Item(5, 4).Value = 1943
$ws.Cells.Item(5, 5).Value = 1992.1
$ws.Cells.Item(5, 6).Value = 1235
$ws.Cells.Item(5, 7).Value = 2.464735706038845
$ws.Cells.Item(5, 8).Value = 2.527020072053521
$ws.Cells.Item(5, 9).Value = 57.32793522267205

# Row 6: BHARTIARTL
$ws.Cells.Item(6, 1).Value = 'INE397D01024'
$ws.Cells.Item(6, 2).Value = 'BHARTIARTL'
$ws.Cells.Item(6, 3).Value = 'Telecom - Cellular & Fixed line services'
$ws.Cells.Item(6, 4).Value = 1594
$ws.Cells.Item(6, 5).Value = 1779
$ws.Cells.Item(6, 6).Value = 1021.35
$ws.Cells.Item(6, 7).Value = 10.3991006183249
$ws.Cells.Item(6, 8).Value = 11.60602258469259
$ws.Cells.Item(6, 9).Value = 56.06794928281196

# Row 7: SHRIRAMFIN
$ws.Cells.Item(7, 1).Value = 'INE721A01013'
$ws.Cells.Item(7, 2).Value = 'SHRIRAMFIN'
$ws.Cells.Item(7, 3).Value = 'Non Banking Financial Company (NBFC)'
$ws.Cells.Item(7, 4).Value = 3055
$ws.Cells.Item(7, 5).Value = 3652.25
$ws.Cells.Item(7, 6).Value = 2029
$ws.Cells.Item(7, 7).Value = 16.35293312341707
$ws.Cells.Item(7, 8).Value = 19.54991816693945
$ws.Cells.Item(7, 9).Value = 50.56678166584525

# Row 8: EICHERMOT
$ws.Cells.Item(8, 1).Value = 'INE066A01021'
$ws.Cells.Item(8, 2).Value = 'EICHERMOT'
$ws.Cells.Item(8, 3).Value = '2/3 Wheelers'
$ws.Cells.Item(8, 4).Value = 5310
$ws.Cells.Item(8, 5).Value = 5385.7
$ws.Cells.Item(8, 6).Value = 3562.45
$ws.Cells.Item(8, 7).Value = 1.405574020090239
$ws.Cells.Item(8, 8).Value = 1.425612052730685
$ws.Cells.Item(8, 9).Value = 49.05472357506773

# Row 9: TECHM
$ws.Cells.Item(9, 1).Value = 'INE669C01036'
$ws.Cells.Item(9, 2).Value = 'TECHM'
$ws.Cells.Item(9, 3).Value = 'Computers - Software & Consulting'
$ws.Cells.Item(9, 4).Value = 1690.5
$ws.Cells.Item(9, 5).Value = 1807.7
$ws.Cells.Item(9, 6).Value = 1162.95
$ws.Cells.Item(9, 7).Value = 6.483376666482277
$ws.Cells.Item(9, 8).Value = 6.932860100561977
$ws.Cells.Item(9, 9).Value = 45.36308525731974

# Row 10: SUNPHARMA
$ws.Cells.Item(10, 1).Value = 'INE044A01036'
$ws.Cells.Item(10, 2).Value = 'SUNPHARMA'
$ws.Cells.Item(10, 3).Value = 'Pharmaceuticals'
$ws.Cells.Item(10, 4).Value = 1848.65
$ws.Cells.Item(10, 5).Value = 1960.35
$ws.Cells.Item(10, 6).Value = 1287
$ws.Cells.Item(10, 7).Value = 5.697962098604837
$ws.Cells.Item(10, 8).Value = 6.042247045141047
$ws.Cells.Item(10, 9).Value = 43.64024864024864

# Row 11: INFY
$ws.Cells.Item(11, 1).Value = 'INE009A01021'
$ws.Cells.Item(11, 2).Value = 'INFY'
$ws.Cells.Item(11, 3).Value = 'Computers - Software & Consulting'
$ws.Cells.Item(11, 4).Value = 1939.1
$ws.Cells.Item(11, 5).Value = 2006.45
$ws.Cells.Item(11, 6).Value = 1358.35
$ws.Cells.Item(11, 7).Value = 3.356674724015063
$ws.Cells.Item(11, 8).Value = 3.47326079108865
$ws.Cells.Item(11, 9).Value = 42.75407663709647

# Row 12: WIPRO
$ws.Cells.Item(12, 1).Value = 'INE075A01022'
$ws.Cells.Item(12, 2).Value = 'WIPRO'
$ws.Cells.Item(12, 3).Value = 'Computers - Software & Consulting'
$ws.Cells.Item(12, 4).Value = 295.2
$ws.Cells.Item(12, 5).Value = 320
$ws.Cells.Item(12, 6).Value = 208.5
$ws.Cells.Item(12, 7).Value = 7.750000000000002
$ws.Cells.Item(12, 8).Value = 8.40108401084012
$ws.Cells.Item(12, 9).Value = 41.58273381294963

# Row 13: POWERGRID
$ws.Cells.Item(13, 1).Value = 'INE752E01010'
$ws.Cells.Item(13, 2).Value = 'POWERGRID'
$ws.Cells.Item(13, 3).Value = 'Power - Transmission'
$ws.Cells.Item(13, 4).Value = 315.8
$ws.Cells.Item(13, 5).Value = 366.25
$ws.Cells.Item(13, 6).Value = 226.05
$ws.Cells.Item(13, 7).Value = 13.77474402730375
$ws.Cells.Item(13, 8).Value = 15.97530082330589
$ws.Cells.Item(13, 9).Value = 39.70360539703606

# Row 14: BAJAJ-AUTO
$ws.Cells.Item(14, 1).Value = 'INE917I01010'
$ws.Cells.Item(14, 2).Value = 'BAJAJ-AUTO'
$ws.Cells.Item(14, 3).Value = '2/3 Wheelers'
$ws.Cells.Item(14, 4).Value = 8965
$ws.Cells.Item(14, 5).Value = 12774
$ws.Cells.Item(14, 6).Value = 6604
$ws.Cells.Item(14, 7).Value = 29.81838108658213
$ws.Cells.Item(14, 8).Value = 42.48745119910764
$ws.Cells.Item(14, 9).Value = 35.75105996365839

# Row 15: BPCL
$ws.Cells.Item(15, 1).Value = 'INE029A01011'
$ws.Cells.Item(15, 2).Value = 'BPCL'
$ws.Cells.Item(15, 3).Value = 'Refineries & Marketing'
$ws.Cells.Item(15, 4).Value = 296.7
$ws.Cells.Item(15, 5).Value = 376
$ws.Cells.Item(15, 6).Value = 222.55
$ws.Cells.Item(15, 7).Value = 21.09042553191489
$ws.Cells.Item(15, 8).Value = 26.72733400741489
$ws.Cells.Item(15, 9).Value = 33.31835542574701

# Row 16: SBIN
$ws.Cells.Item(16, 1).Value = 'INE062A01020'
$ws.Cells.Item(16, 2).Value = 'SBIN'
$ws.Cells.Item(16, 3).Value = 'Public Sector Bank'
$ws.Cells.Item(16, 4).Value = 793.3
$ws.Cells.Item(16, 5).Value = 912
$ws.Cells.Item(16, 6).Value = 600.65
$ws.Cells.Item(16, 7).Value = 13.01535087719299
$ws.Cells.Item(16, 8).Value = 14.96281356359512
$ws.Cells.Item(16, 9).Value = 32.07358694747357

# Row 17: ICICIBANK
$ws.Cells.Item(17, 1).Value = 'INE090A01021'
$ws.Cells.Item(17, 2).Value = 'ICICIBANK'
$ws.Cells.Item(17, 3).Value = 'Private Sector Bank'
$ws.Cells.Item(17, 4).Value = 1265.6
$ws.Cells.Item(17, 5).Value = 1362.35
$ws.Cells.Item(17, 6).Value = 970.15
$ws.Cells.Item(17, 7).Value = 7.101699269644368
$ws.Cells.Item(17, 8).Value = 7.644595448798985
$ws.Cells.Item(17, 9).Value = 30.45405349688193
